$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.15%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'37.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'7.03%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-3.40%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07840"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.51%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.211"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-4.51%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'8.014"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.29%"
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'1.13%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9160"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-1.33%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.09675"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-4.87%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1883"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'3.85%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08594"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.18%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03569"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'2.40%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09965"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.55%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001480"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.62%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005671"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.56%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.460"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.35%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.398"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'13.82%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3463"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.62%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1319"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.47%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.775"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'5.31%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'-1.61%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'0.09%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001233"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.00%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004781"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'7.82%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'7.88%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'39.77%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01787"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'1.93%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04741"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.53%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.008126"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'5.91%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1392"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-1.01%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007681"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'9.08%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002142"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-5.98%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009943"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'6.65%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006190"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'3.41%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.32%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'7.998"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'192.73%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-0.25%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.32%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.32%"
$ws.Range("E51").Style = "Normal"
